# Refresh the cryptos price/volume snapshot (GitHub Actions data pull),
# plus the OKB / dogwifhat row-order swap (rows 45-46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $cell = $ws.Range($rangeAddress)
    # Force text storage so numeric-looking strings (e.g. "594.28")
    # are not silently re-typed as numbers by Excel, then drop the
    # now-unneeded explicit format so the cell keeps the workbook's
    # original (unstyled) look.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue 'D2' '66.788.34'
Set-TextValue 'E2' '  +0.37%  '
Set-TextValue 'D3' '3.493.55'
Set-TextValue 'E3' '  +0.10%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '594.28'
Set-TextValue 'E5' '  +0.63%  '
Set-TextValue 'D6' '172.32'
Set-TextValue 'E6' '  +2.49%  '
Set-TextValue 'D7' '0.999'
Set-TextValue 'E7' '  +0.01%  '
Set-TextValue 'E8' '  -2.18%  '
Set-TextValue 'E9' '  +3.00%  '
Set-TextValue 'E10' '  -2.69%  '
Set-TextValue 'E11' '  -0.14%  '
Set-TextValue 'D12' '4.095.92'
Set-TextValue 'E12' '  +0.09%  '
Set-TextValue 'E13' '  +0.20%  '
Set-TextValue 'D14' '29.21'
Set-TextValue 'E14' '  +4.36%  '
Set-TextValue 'D15' '66.762.25'
Set-TextValue 'D16' '0.0000178'
Set-TextValue 'E16' '  -0.24%  '
Set-TextValue 'D17' '3.496.32'
Set-TextValue 'E17' '  +0.67%  '
Set-TextValue 'D18' '6.25'
Set-TextValue 'E18' '  -0.17%  '
Set-TextValue 'D19' '14.29'
Set-TextValue 'E19' '  +2.53%  '
Set-TextValue 'D20' '389.98'
Set-TextValue 'E20' '  -0.09%  '
Set-TextValue 'D21' '7.90'
Set-TextValue 'E21' '  +0.28%  '
Set-TextValue 'D22' '73.30'
Set-TextValue 'E22' '  +0.76%  '
Set-TextValue 'E23' '  +0.01%  '
Set-TextValue 'D24' '0.533'
Set-TextValue 'E24' '  +0.50%  '
Set-TextValue 'E25' '  -0.66%  '
Set-TextValue 'D26' '0.0000121'
Set-TextValue 'E26' '  -0.52%  '
Set-TextValue 'E27' '  -0.82%  '
Set-TextValue 'E28' '  +0.01%  '
Set-TextValue 'E29' '  -0.55%  '
Set-TextValue 'E30' '  -3.48%  '
Set-TextValue 'E31' '  -1.85%  '
Set-TextValue 'E32' '  +0.10%  '
Set-TextValue 'D33' '23.55'
Set-TextValue 'E33' '  -0.28%  '
Set-TextValue 'D34' '7.34'
Set-TextValue 'E34' '  +0.46%  '
Set-TextValue 'D35' '1.60'
Set-TextValue 'E35' '  +0.96%  '
Set-TextValue 'D36' '163.60'
Set-TextValue 'E36' '  +0.50%  '
Set-TextValue 'D37' '0.875'
Set-TextValue 'E37' '  -2.22%  '
Set-TextValue 'D38' '1.90'
Set-TextValue 'E38' '  -0.42%  '
Set-TextValue 'D39' '6.81'
Set-TextValue 'E39' '  +0.25%  '
Set-TextValue 'D40' '4.62'
Set-TextValue 'E40' '  -0.11%  '
Set-TextValue 'D41' '2.817.64'
Set-TextValue 'E41' '  +1.78%  '
Set-TextValue 'E42' '  +1.61%  '
Set-TextValue 'D43' '0.0727'
Set-TextValue 'E43' '  -1.20%  '
Set-TextValue 'D44' '25.90'
Set-TextValue 'E44' '  -0.94%  '
Set-TextValue 'B45' 'dogwifhat'
Set-TextValue 'C45' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D45' '2.55'
Set-TextValue 'E45' '  +0.16%  '
Set-TextValue 'B46' 'OKB'
Set-TextValue 'C46' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D46' '42.40'
Set-TextValue 'E46' '  -0.91%  '
Set-TextValue 'E47' '  -2.80%  '
Set-TextValue 'D48' '337.92'
Set-TextValue 'E48' '  -1.07%  '
Set-TextValue 'E49' '  -0.52%  '
Set-TextValue 'D50' '33.74'
Set-TextValue 'E50' '  +1.12%  '
Set-TextValue 'D51' '6.40'
Set-TextValue 'E51' '  -0.90%  '
